$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$updates = @{
    2  = 5
    7  = 5
    16 = -1
    20 = -6
    28 = 1
    34 = -4
    36 = 2
    42 = 1
    44 = 8
    48 = 0
    61 = -1
    63 = -2
    66 = -2
    69 = -1
    71 = -2
    74 = -3
    76 = -6
    78 = 3
    81 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
